$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44508
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("P2").Value = 667

# Row 3
$ws.Range("D3").Value = 44518
$ws.Range("J3").Value = 50

# Row 4
$ws.Range("D4").Value = 44749
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 20000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 20000
$ws.Range("P4").Value = 1333

# Row 5
$ws.Range("D5").Value = 44525
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 8000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 8000
$ws.Range("P5").Value = 533
